# "Generate Report for Handoff"
#
# The localization-status report is regenerated: two files ("low" priority
# items) get re-handed-off, which bumps their "Latest Handoff Datetime"
# stamps and flips their Priority from "low" to "ht" in both language
# sheets (zh-cn, de-de). The Overview sheet's "Latest HO Xliff Generate
# Date" column (shared with the de-de sheet's handoff datetime, since they
# were generated together) is updated to match.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for the 4 rows that were
# regenerated (rows 4-7).
$ws1.Range("G4:G7").Value = "2016-08-24 04:29:36"

# zh-cn: Priority low -> ht, and Latest Handoff Datetime refreshed.
$ws2.Range("E4:E7").Value = "ht"
$ws2.Range("H4:H7").Value = "2016-08-24 04:29:31"

# de-de: Priority low -> ht, and Latest Handoff Datetime refreshed
# (matches the Overview generate date above).
$ws3.Range("E4:E7").Value = "ht"
$ws3.Range("H4:H7").Value = "2016-08-24 04:29:36"
